# Applies arithmetic-expression replacements across the table cells.
$d = $word.ActiveDocument

$replacements = @(
    @("9-5=", "31+31="),
    @("88-71=", "3+36="),
    @("83+9=", "1+23="),
    @("93-13=", "35-20="),
    @("32-11=", "86-32="),
    @("84-18=", "92-51="),
    @("15-5=", "46+0="),
    @("17+29=", "54+42="),
    @("47-20=", "99-72="),
    @("9+39=", "24-13="),
    @("98-24=", "17+54="),
    @("87-34=", "52+23="),
    @("43+52=", "25-11="),
    @("28+68=", "97-39="),
    @("65+11=", "21+46="),
    @("27-17=", "27+26="),
    @("74+16=", "76-52="),
    @("83+14=", "66-62="),
    @("33+1=", "65-57="),
    @("61+36=", "57-47="),
    @("4+92=", "3+29="),
    @("24-21=", "23+74="),
    @("31+62=", "86-27="),
    @("4+95=", "76+1="),
    @("28+69=", "11+38="),
    @("87-7=", "29+11="),
    @("96-23=", "38+61="),
    @("80-21=", "40+4="),
    @("22+20=", "25+0="),
    @("93-82=", "4+25="),
    @("97-6=", "32-28="),
    @("37+6=", "64-8="),
    @("67-51=", "48-14="),
    @("73+26=", "21-16="),
    @("67-5=", "15+69="),
    @("75-32=", "68+23="),
    @("56-44=", "86-77="),
    @("52+25=", "14+10="),
    @("25+66=", "23+52="),
    @("62-6=", "72+4="),
    @("66-34=", "53-6="),
    @("60-12=", "71-5="),
    @("51-36=", "65+15="),
    @("99-9=", "68+27="),
    @("47+11=", "60+13="),
    @("5+57=", "45-15="),
    @("18+25=", "11+24="),
    @("99-3=", "75+23="),
    @("65-51=", "55+8="),
    @("76-72=", "90-75="),
    @("52+21=", "10+61="),
    @("90-72=", "50+7="),
    @("83-65=", "16-3="),
    @("90-45=", "74-32="),
    @("64-57=", "58+30="),
    @("42+12=", "99-37="),
    @("30-2=", "62+31="),
    @("35+11=", "14+68="),
    @("64-13=", "40+48="),
    @("72-16=", "62-16="),
    @("47+48=", "11+15="),
    @("48-17=", "84-14="),
    @("68+21=", "9+33="),
    @("92-83=", "4+9="),
    @("88-79=", "23-17="),
    @("61+15=", "11+74="),
    @("64-16=", "49+12="),
    @("10+57=", "20+1="),
    @("3+45=", "70-9="),
    @("17-3=", "45-36="),
    @("81-55=", "99-28="),
    @("61+27=", "94-88="),
    @("36-17=", "73+19="),
    @("33+51=", "87-75="),
    @("33-11=", "13+51="),
    @("4+78=", "78-71="),
    @("38+21=", "86-18="),
    @("84+10=", "88-67="),
    @("68-54=", "8+39="),
    @("55-24=", "96-62="),
    @("90-14=", "13+30="),
    @("42+48=", "14+79="),
    @("64+34=", "19+8="),
    @("44-40=", "64+4="),
    @("6+37=", "11+50="),
    @("80-37=", "36+12="),
    @("59-18=", "35+31="),
    @("77-52=", "34+30="),
    @("27-16=", "53+16="),
    @("53+3=", "43+39="),
    @("95-27=", "62+21="),
    @("61+17=", "45+7="),
    @("64+28=", "82-76="),
    @("26-10=", "95-5="),
    @("7+92=", "28-0="),
    @("31+51=", "66-20="),
    @("73-20=", "32+2="),
    @("5+5=", "96-35="),
    @("17+47=", "18+10="),
    @("42+35=", "15+35=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Applied $($replacements.Count) replacements"
